$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 92 from 45186 to 45188
for ($row = 2; $row -le 92; $row++) {
    $ws.Cells.Item($row, 3).Value = 45188
}
